# Beilage_Verfuegung_per_Kreditor.xlsx
# On the "AELCOM AG" sheet, remove the trailing "Begründungen (NA15)" block
# (rows 26-28: the section header, its column headings, and the single
# "963040.0 / Anderes" data row, along with the merged cells B27:F27 and
# B28:F28 that belonged to it). Deleting the entire rows shifts nothing up
# from below (they were the last rows on the sheet) and shrinks the sheet's
# used range from A2:H28 down to A2:H25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AELCOM AG")

$ws.Rows("26:28").EntireRow.Delete()
